$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data row (row 1: A1, B1, C1) entirely.
$ws.Rows.Item(1).Delete()

# Write the new data row (row 3): a number, two date-like numbers, a string, and a boolean.
$ws.Range("A3").Value = 1.1
$ws.Range("B3").Value = 42894.36889262732
$ws.Range("C3").Value = 42894.36889262732
$ws.Range("D3").Value = "a string"
$ws.Range("E3").Value = $true
